$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("B24").Value = 6206197
$ws.Range("E24").Value = 'Chungbuk Cheongju'
$ws.Range("F24").Value = 'Seongnam FC'
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 'D'
$ws.Range("L24").Value = 2.7
$ws.Range("M24").Value = 3.1
$ws.Range("N24").Value = 2.5
$ws.Range("O24").Value = 2.7
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = 2.55
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = 1.925
$ws.Range("T24").Value = 1.875
$ws.Range("U24").Value = 2.25
$ws.Range("V24").Value = 1.85
$ws.Range("W24").Value = 1.95
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = 2
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0
$ws.Range("AB24").Value = 0
$ws.Range("AC24").Value = -1
$ws.Range("AD24").Value = 0.95

# Row 25
$ws.Range("B25").Value = 6204317
$ws.Range("E25").Value = 'Jeonnam Dragons'
$ws.Range("F25").Value = 'Seoul ELand FC'
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 2
$ws.Range("K25").Value = 'D'
$ws.Range("L25").Value = 2.4
$ws.Range("M25").Value = 3.3
$ws.Range("N25").Value = 2.625
$ws.Range("O25").Value = 2.15
$ws.Range("P25").Value = 3.4
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = -0.25
$ws.Range("S25").Value = 1.9
$ws.Range("T25").Value = 1.9
$ws.Range("U25").Value = 2.5
$ws.Range("V25").Value = 1.975
$ws.Range("W25").Value = 1.825
$ws.Range("X25").Value = -1
$ws.Range("Y25").Value = 2.4
$ws.Range("Z25").Value = -1
$ws.Range("AA25").Value = -0.5
$ws.Range("AB25").Value = 0.45
$ws.Range("AC25").Value = 0.9750000000000001
$ws.Range("AD25").Value = -1

# Row 41
$ws.Range("B41").Value = 6206211
$ws.Range("E41").Value = 'Seongnam FC'
$ws.Range("F41").Value = 'Gyeongnam FC'
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 1
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 'D'
$ws.Range("L41").Value = 2.55
$ws.Range("M41").Value = 2.875
$ws.Range("N41").Value = 2.7
$ws.Range("O41").Value = 2.7
$ws.Range("P41").Value = 2.9
$ws.Range("Q41").Value = 2.5
$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 1.975
$ws.Range("T41").Value = 1.825
$ws.Range("U41").Value = 2.5
$ws.Range("V41").Value = 1.95
$ws.Range("W41").Value = 1.85
$ws.Range("X41").Value = -1
$ws.Range("Y41").Value = 1.9
$ws.Range("Z41").Value = -1
$ws.Range("AA41").Value = 0
$ws.Range("AB41").Value = 0
$ws.Range("AC41").Value = -1
$ws.Range("AD41").Value = 0.8500000000000001

# Row 42
$ws.Range("B42").Value = 6204320
$ws.Range("E42").Value = 'Jeonnam Dragons'
$ws.Range("F42").Value = 'Ansan Greeners FC'
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 2
$ws.Range("I42").Value = 3
$ws.Range("J42").Value = 1
$ws.Range("K42").Value = 'H'
$ws.Range("L42").Value = 1.666
$ws.Range("M42").Value = 3.75
$ws.Range("N42").Value = 4.333
$ws.Range("O42").Value = 1.571
$ws.Range("P42").Value = 4
$ws.Range("Q42").Value = 4.5
$ws.Range("R42").Value = -1
$ws.Range("S42").Value = 2.025
$ws.Range("T42").Value = 1.775
$ws.Range("U42").Value = 2.75
$ws.Range("V42").Value = 1.8
$ws.Range("W42").Value = 2
$ws.Range("X42").Value = 0.571
$ws.Range("Y42").Value = -1
$ws.Range("Z42").Value = -1
$ws.Range("AA42").Value = 1.025
$ws.Range("AB42").Value = -1
$ws.Range("AC42").Value = 0.8
$ws.Range("AD42").Value = -1

# Row 80
$ws.Range("B80").Value = 6206245
$ws.Range("E80").Value = 'Ansan Greeners FC'
$ws.Range("F80").Value = 'Bucheon'
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2
$ws.Range("K80").Value = 'A'
$ws.Range("L80").Value = 3.1
$ws.Range("M80").Value = 3.3
$ws.Range("N80").Value = 2.15
$ws.Range("O80").Value = 3.2
$ws.Range("P80").Value = 3.1
$ws.Range("Q80").Value = 2.2
$ws.Range("R80").Value = 0.25
$ws.Range("S80").Value = 1.875
$ws.Range("T80").Value = 1.925
$ws.Range("U80").Value = 2.25
$ws.Range("V80").Value = 1.975
$ws.Range("W80").Value = 1.825
$ws.Range("X80").Value = -1
$ws.Range("Y80").Value = -1
$ws.Range("Z80").Value = 1.2
$ws.Range("AA80").Value = -1
$ws.Range("AB80").Value = 0.925
$ws.Range("AC80").Value = 0.9750000000000001
$ws.Range("AD80").Value = -1

# Row 81
$ws.Range("B81").Value = 6206246
$ws.Range("E81").Value = 'Busan I Park'
$ws.Range("F81").Value = 'Gimcheon Sangmu FC'
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 'H'
$ws.Range("L81").Value = 2.9
$ws.Range("M81").Value = 3.4
$ws.Range("N81").Value = 2.2
$ws.Range("O81").Value = 2.9
$ws.Range("P81").Value = 3.4
$ws.Range("Q81").Value = 2.25
$ws.Range("R81").Value = 0.25
$ws.Range("S81").Value = 1.85
$ws.Range("T81").Value = 1.95
$ws.Range("U81").Value = 2.5
$ws.Range("V81").Value = 1.875
$ws.Range("W81").Value = 1.925
$ws.Range("X81").Value = 1.9
$ws.Range("Y81").Value = -1
$ws.Range("Z81").Value = -1
$ws.Range("AA81").Value = 0.8500000000000001
$ws.Range("AB81").Value = -1
$ws.Range("AC81").Value = -1
$ws.Range("AD81").Value = 0.925

# Row 83
$ws.Range("B83").Value = 6206247
$ws.Range("E83").Value = 'Chungbuk Cheongju'
$ws.Range("F83").Value = 'FC Anyang'
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 1
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 'H'
$ws.Range("L83").Value = 3.5
$ws.Range("M83").Value = 3.3
$ws.Range("N83").Value = 2
$ws.Range("O83").Value = 3.5
$ws.Range("P83").Value = 3.2
$ws.Range("Q83").Value = 2
$ws.Range("R83").Value = 0.25
$ws.Range("S83").Value = 1.975
$ws.Range("T83").Value = 1.825
$ws.Range("U83").Value = 2.25
$ws.Range("V83").Value = 2.025
$ws.Range("W83").Value = 1.775
$ws.Range("X83").Value = 2.5
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.9750000000000001
$ws.Range("AB83").Value = -1
$ws.Range("AC83").Value = 1.025
$ws.Range("AD83").Value = -1

# Row 84
$ws.Range("B84").Value = 6206248
$ws.Range("E84").Value = 'Seongnam FC'
$ws.Range("F84").Value = 'Gimpo FC'
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 4
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 'A'
$ws.Range("L84").Value = 2.375
$ws.Range("M84").Value = 3.2
$ws.Range("N84").Value = 2.875
$ws.Range("O84").Value = 2
$ws.Range("P84").Value = 3.3
$ws.Range("Q84").Value = 3.6
$ws.Range("R84").Value = -0.25
$ws.Range("S84").Value = 1.75
$ws.Range("T84").Value = 2.05
$ws.Range("U84").Value = 2
$ws.Range("V84").Value = 1.85
$ws.Range("W84").Value = 1.95
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 2.6
$ws.Range("AA84").Value = -1
$ws.Range("AB84").Value = 1.05
$ws.Range("AC84").Value = 0.8500000000000001
$ws.Range("AD84").Value = -1

# Row 91
$ws.Range("B91").Value = 6206255
$ws.Range("E91").Value = 'Seoul ELand FC'
$ws.Range("F91").Value = 'Cheonan City'
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 3
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1
$ws.Range("K91").Value = 'A'
$ws.Range("L91").Value = 1.666
$ws.Range("M91").Value = 3.75
$ws.Range("N91").Value = 4.5
$ws.Range("O91").Value = 1.833
$ws.Range("P91").Value = 3.5
$ws.Range("Q91").Value = 3.8
$ws.Range("R91").Value = -0.5
$ws.Range("S91").Value = 1.825
$ws.Range("T91").Value = 1.975
$ws.Range("U91").Value = 2.25
$ws.Range("V91").Value = 1.825
$ws.Range("W91").Value = 1.975
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = 2.8
$ws.Range("AA91").Value = -1
$ws.Range("AB91").Value = 0.9750000000000001
$ws.Range("AC91").Value = 0.825
$ws.Range("AD91").Value = -1

# Row 92
$ws.Range("B92").Value = 6206254
$ws.Range("E92").Value = 'Ansan Greeners FC'
$ws.Range("F92").Value = 'Chungbuk Cheongju'
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 1
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 'A'
$ws.Range("L92").Value = 3.1
$ws.Range("M92").Value = 3.4
$ws.Range("N92").Value = 2.1
$ws.Range("O92").Value = 3.1
$ws.Range("P92").Value = 3.3
$ws.Range("Q92").Value = 2.15
$ws.Range("R92").Value = 0.25
$ws.Range("S92").Value = 1.875
$ws.Range("T92").Value = 1.925
$ws.Range("U92").Value = 2.25
$ws.Range("V92").Value = 1.8
$ws.Range("W92").Value = 2
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = -1
$ws.Range("Z92").Value = 1.15
$ws.Range("AA92").Value = -1
$ws.Range("AB92").Value = 0.925
$ws.Range("AC92").Value = -1
$ws.Range("AD92").Value = 1

# Row 148
$ws.Range("B148").Value = 7737361
$ws.Range("E148").Value = 'FC Anyang'
$ws.Range("F148").Value = 'Seongnam FC'
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 0
$ws.Range("I148").Value = 2
$ws.Range("J148").Value = 0
$ws.Range("K148").Value = 'H'
$ws.Range("L148").Value = 1.833
$ws.Range("M148").Value = 3.5
$ws.Range("N148").Value = 4.333
$ws.Range("O148").Value = 2.15
$ws.Range("P148").Value = 3.4
$ws.Range("Q148").Value = 3.3
$ws.Range("R148").Value = -0.25
$ws.Range("S148").Value = 1.85
$ws.Range("T148").Value = 1.95
$ws.Range("U148").Value = 2.5
$ws.Range("V148").Value = 1.925
$ws.Range("W148").Value = 1.875
$ws.Range("X148").Value = 1.15
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = -1
$ws.Range("AA148").Value = 0.8500000000000001
$ws.Range("AB148").Value = -1
$ws.Range("AC148").Value = -1
$ws.Range("AD148").Value = 0.875

# Row 149
$ws.Range("B149").Value = 7738678
$ws.Range("E149").Value = 'Ansan Greeners FC'
$ws.Range("F149").Value = 'Gyeongnam FC'
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 2
$ws.Range("I149").Value = 1
$ws.Range("J149").Value = 1
$ws.Range("K149").Value = 'A'
$ws.Range("L149").Value = 4
$ws.Range("M149").Value = 3.6
$ws.Range("N149").Value = 1.85
$ws.Range("O149").Value = 4.333
$ws.Range("P149").Value = 3.5
$ws.Range("Q149").Value = 1.833
$ws.Range("R149").Value = 0.5
$ws.Range("S149").Value = 1.95
$ws.Range("T149").Value = 1.85
$ws.Range("U149").Value = 2.5
$ws.Range("V149").Value = 1.9
$ws.Range("W149").Value = 1.9
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 0.833
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = 0.8500000000000001
$ws.Range("AC149").Value = 0.8999999999999999
$ws.Range("AD149").Value = -1

# Row 161
$ws.Range("B161").Value = 7738657
$ws.Range("E161").Value = 'Bucheon'
$ws.Range("F161").Value = 'Seoul ELand FC'
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 0
$ws.Range("I161").Value = 0
$ws.Range("J161").Value = 0
$ws.Range("K161").Value = 'H'
$ws.Range("L161").Value = 2.5
$ws.Range("M161").Value = 3.2
$ws.Range("N161").Value = 2.5
$ws.Range("O161").Value = 2.45
$ws.Range("P161").Value = 3.2
$ws.Range("Q161").Value = 2.55
$ws.Range("R161").Value = 0
$ws.Range("S161").Value = 1.875
$ws.Range("T161").Value = 1.925
$ws.Range("U161").Value = 2.25
$ws.Range("V161").Value = 1.825
$ws.Range("W161").Value = 1.975
$ws.Range("X161").Value = 1.45
$ws.Range("Y161").Value = -1
$ws.Range("Z161").Value = -1
$ws.Range("AA161").Value = 0.875
$ws.Range("AB161").Value = -1
$ws.Range("AC161").Value = -1
$ws.Range("AD161").Value = 0.9750000000000001

# Row 162
$ws.Range("B162").Value = 7738682
$ws.Range("E162").Value = 'Gyeongnam FC'
$ws.Range("F162").Value = 'Chungnam Asan FC'
$ws.Range("G162").Value = 1
$ws.Range("H162").Value = 2
$ws.Range("I162").Value = 1
$ws.Range("J162").Value = 2
$ws.Range("K162").Value = 'A'
$ws.Range("L162").Value = 2.25
$ws.Range("M162").Value = 3.25
$ws.Range("N162").Value = 2.75
$ws.Range("O162").Value = 2
$ws.Range("P162").Value = 3.3
$ws.Range("Q162").Value = 3.2
$ws.Range("R162").Value = -0.25
$ws.Range("S162").Value = 1.8
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 2.25
$ws.Range("V162").Value = 1.875
$ws.Range("W162").Value = 1.925
$ws.Range("X162").Value = -1
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = 2.2
$ws.Range("AA162").Value = -1
$ws.Range("AB162").Value = 1
$ws.Range("AC162").Value = 0.875
$ws.Range("AD162").Value = -1

# Row 163
$ws.Range("B163").Value = 7738683
$ws.Range("E163").Value = 'Chungbuk Cheongju'
$ws.Range("F163").Value = 'FC Anyang'
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 1
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = 0
$ws.Range("K163").Value = 'D'
$ws.Range("L163").Value = 2.3
$ws.Range("M163").Value = 3.2
$ws.Range("N163").Value = 2.7
$ws.Range("O163").Value = 2.75
$ws.Range("P163").Value = 3.2
$ws.Range("Q163").Value = 2.3
$ws.Range("R163").Value = 0.25
$ws.Range("S163").Value = 1.75
$ws.Range("T163").Value = 2.05
$ws.Range("U163").Value = 2.25
$ws.Range("V163").Value = 1.925
$ws.Range("W163").Value = 1.875
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = 2.2
$ws.Range("Z163").Value = -1
$ws.Range("AA163").Value = 0.375
$ws.Range("AB163").Value = -0.5
$ws.Range("AC163").Value = -0.5
$ws.Range("AD163").Value = 0.4375

# Row 164
$ws.Range("B164").Value = 7737346
$ws.Range("E164").Value = 'Busan I Park'
$ws.Range("F164").Value = 'Gimpo FC'
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1
$ws.Range("I164").Value = 0
$ws.Range("J164").Value = 0
$ws.Range("K164").Value = 'A'
$ws.Range("L164").Value = 1.8
$ws.Range("M164").Value = 3.25
$ws.Range("N164").Value = 4
$ws.Range("O164").Value = 1.7
$ws.Range("P164").Value = 3.3
$ws.Range("Q164").Value = 4.5
$ws.Range("R164").Value = -0.75
$ws.Range("S164").Value = 2
$ws.Range("T164").Value = 1.8
$ws.Range("U164").Value = 2.25
$ws.Range("V164").Value = 1.9
$ws.Range("W164").Value = 1.9
$ws.Range("X164").Value = -1
$ws.Range("Y164").Value = -1
$ws.Range("Z164").Value = 3.5
$ws.Range("AA164").Value = -1
$ws.Range("AB164").Value = 0.8
$ws.Range("AC164").Value = -1
$ws.Range("AD164").Value = 0.8999999999999999

# Row 175
$ws.Range("B175").Value = 7738661
$ws.Range("E175").Value = 'Cheonan City'
$ws.Range("F175").Value = 'Gyeongnam FC'
$ws.Range("G175").Value = 2
$ws.Range("H175").Value = 2
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 1
$ws.Range("K175").Value = 'D'
$ws.Range("L175").Value = 4
$ws.Range("M175").Value = 3.4
$ws.Range("N175").Value = 1.833
$ws.Range("O175").Value = 3.6
$ws.Range("P175").Value = 3.3
$ws.Range("Q175").Value = 1.95
$ws.Range("R175").Value = 0.5
$ws.Range("S175").Value = 1.8
$ws.Range("T175").Value = 2
$ws.Range("U175").Value = 2.5
$ws.Range("V175").Value = 2.025
$ws.Range("W175").Value = 1.775
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = 2.3
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = 0.8
$ws.Range("AB175").Value = -1
$ws.Range("AC175").Value = 1.025
$ws.Range("AD175").Value = -1

# Row 176
$ws.Range("B176").Value = 7737347
$ws.Range("E176").Value = 'Busan I Park'
$ws.Range("F176").Value = 'Jeonnam Dragons'
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 1
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 0
$ws.Range("K176").Value = 'A'
$ws.Range("L176").Value = 1.909
$ws.Range("M176").Value = 3.25
$ws.Range("N176").Value = 3.75
$ws.Range("O176").Value = 1.909
$ws.Range("P176").Value = 3.25
$ws.Range("Q176").Value = 3.8
$ws.Range("R176").Value = -0.5
$ws.Range("S176").Value = 1.925
$ws.Range("T176").Value = 1.875
$ws.Range("U176").Value = 2.25
$ws.Range("V176").Value = 1.95
$ws.Range("W176").Value = 1.85
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = 2.8
$ws.Range("AA176").Value = -1
$ws.Range("AB176").Value = 0.875
$ws.Range("AC176").Value = -1
$ws.Range("AD176").Value = 0.8500000000000001

# Row 182
$ws.Range("B182").Value = 7737348
$ws.Range("E182").Value = 'Cheonan City'
$ws.Range("F182").Value = 'Busan I Park'
$ws.Range("G182").Value = 2
$ws.Range("H182").Value = 4
$ws.Range("I182").Value = 1
$ws.Range("J182").Value = 1
$ws.Range("K182").Value = 'A'
$ws.Range("L182").Value = 3
$ws.Range("M182").Value = 3.2
$ws.Range("N182").Value = 2.1
$ws.Range("O182").Value = 3.5
$ws.Range("P182").Value = 3.3
$ws.Range("Q182").Value = 1.909
$ws.Range("R182").Value = 0.5
$ws.Range("S182").Value = 1.85
$ws.Range("T182").Value = 1.95
$ws.Range("U182").Value = 2.25
$ws.Range("V182").Value = 1.825
$ws.Range("W182").Value = 1.975
$ws.Range("X182").Value = -1
$ws.Range("Y182").Value = -1
$ws.Range("Z182").Value = 0.909
$ws.Range("AA182").Value = -1
$ws.Range("AB182").Value = 0.95
$ws.Range("AC182").Value = 0.825
$ws.Range("AD182").Value = -1

# Row 183
$ws.Range("B183").Value = 7738689
$ws.Range("E183").Value = 'Ansan Greeners FC'
$ws.Range("F183").Value = 'Chungnam Asan FC'
$ws.Range("G183").Value = 1
$ws.Range("H183").Value = 0
$ws.Range("I183").Value = 0
$ws.Range("J183").Value = 0
$ws.Range("K183").Value = 'H'
$ws.Range("L183").Value = 3
$ws.Range("M183").Value = 3
$ws.Range("N183").Value = 2.25
$ws.Range("O183").Value = 4
$ws.Range("P183").Value = 3.1
$ws.Range("Q183").Value = 1.909
$ws.Range("R183").Value = 0.5
$ws.Range("S183").Value = 1.85
$ws.Range("T183").Value = 1.95
$ws.Range("U183").Value = 2.25
$ws.Range("V183").Value = 1.975
$ws.Range("W183").Value = 1.825
$ws.Range("X183").Value = 3
$ws.Range("Y183").Value = -1
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0.8500000000000001
$ws.Range("AB183").Value = -1
$ws.Range("AC183").Value = -1
$ws.Range("AD183").Value = 0.825

# Row 192
$ws.Range("B192").Value = 7739336
$ws.Range("E192").Value = 'FC Anyang'
$ws.Range("F192").Value = 'Suwon Bluewings'
$ws.Range("G192").Value = 1
$ws.Range("H192").Value = 3
$ws.Range("I192").Value = 0
$ws.Range("J192").Value = 2
$ws.Range("K192").Value = 'A'
$ws.Range("L192").Value = 2.7
$ws.Range("M192").Value = 3.3
$ws.Range("N192").Value = 2.375
$ws.Range("O192").Value = 2.75
$ws.Range("P192").Value = 3.25
$ws.Range("Q192").Value = 2.375
$ws.Range("R192").Value = 0
$ws.Range("S192").Value = 2.05
$ws.Range("T192").Value = 1.75
$ws.Range("U192").Value = 2.25
$ws.Range("V192").Value = 1.775
$ws.Range("W192").Value = 2.025
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = -1
$ws.Range("Z192").Value = 1.375
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = 0.75
$ws.Range("AC192").Value = 0.7749999999999999
$ws.Range("AD192").Value = -1

# Row 193
$ws.Range("B193").Value = 7738692
$ws.Range("E193").Value = 'Ansan Greeners FC'
$ws.Range("F193").Value = 'Gimpo FC'
$ws.Range("G193").Value = 1
$ws.Range("H193").Value = 2
$ws.Range("I193").Value = 0
$ws.Range("J193").Value = 2
$ws.Range("K193").Value = 'A'
$ws.Range("L193").Value = 3.3
$ws.Range("M193").Value = 3.3
$ws.Range("N193").Value = 2.05
$ws.Range("O193").Value = 3.1
$ws.Range("P193").Value = 3.1
$ws.Range("Q193").Value = 2.2
$ws.Range("R193").Value = 0.25
$ws.Range("S193").Value = 1.85
$ws.Range("T193").Value = 1.95
$ws.Range("U193").Value = 2
$ws.Range("V193").Value = 1.75
$ws.Range("W193").Value = 2.05
$ws.Range("X193").Value = -1
$ws.Range("Y193").Value = -1
$ws.Range("Z193").Value = 1.2
$ws.Range("AA193").Value = -1
$ws.Range("AB193").Value = 0.95
$ws.Range("AC193").Value = 0.75
$ws.Range("AD193").Value = -1

# Row 211
$ws.Range("B211").Value = 7737369
$ws.Range("E211").Value = 'Gyeongnam FC'
$ws.Range("F211").Value = 'Seongnam FC'
$ws.Range("G211").Value = 1
$ws.Range("H211").Value = 2
$ws.Range("I211").Value = 1
$ws.Range("J211").Value = 0
$ws.Range("K211").Value = 'A'
$ws.Range("L211").Value = 2.05
$ws.Range("M211").Value = 3.4
$ws.Range("N211").Value = 3.3
$ws.Range("O211").Value = 2.1
$ws.Range("P211").Value = 3.3
$ws.Range("Q211").Value = 3.1
$ws.Range("R211").Value = -0.25
$ws.Range("S211").Value = 1.9
$ws.Range("T211").Value = 1.9
$ws.Range("U211").Value = 2.5
$ws.Range("V211").Value = 1.925
$ws.Range("W211").Value = 1.875
$ws.Range("X211").Value = -1
$ws.Range("Y211").Value = -1
$ws.Range("Z211").Value = 2.1
$ws.Range("AA211").Value = -1
$ws.Range("AB211").Value = 0.8999999999999999
$ws.Range("AC211").Value = 0.925
$ws.Range("AD211").Value = -1

# Row 212
$ws.Range("B212").Value = 7738699
$ws.Range("E212").Value = 'FC Anyang'
$ws.Range("F212").Value = 'Gimpo FC'
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0
$ws.Range("I212").Value = 0
$ws.Range("J212").Value = 0
$ws.Range("K212").Value = 'D'
$ws.Range("L212").Value = 1.75
$ws.Range("M212").Value = 3.4
$ws.Range("N212").Value = 4.333
$ws.Range("O212").Value = 1.833
$ws.Range("P212").Value = 3.3
$ws.Range("Q212").Value = 4
$ws.Range("R212").Value = -0.5
$ws.Range("S212").Value = 1.85
$ws.Range("T212").Value = 1.95
$ws.Range("U212").Value = 2.25
$ws.Range("V212").Value = 1.85
$ws.Range("W212").Value = 1.95
$ws.Range("X212").Value = -1
$ws.Range("Y212").Value = 2.3
$ws.Range("Z212").Value = -1
$ws.Range("AA212").Value = -1
$ws.Range("AB212").Value = 0.95
$ws.Range("AC212").Value = -1
$ws.Range("AD212").Value = 0.95

# Row 218
$ws.Range("B218").Value = 7737353
$ws.Range("E218").Value = 'Chungbuk Cheongju'
$ws.Range("F218").Value = 'Busan I Park'
$ws.Range("G218").Value = 0
$ws.Range("H218").Value = 0
$ws.Range("I218").Value = 0
$ws.Range("J218").Value = 0
$ws.Range("K218").Value = 'D'
$ws.Range("L218").Value = 3.4
$ws.Range("M218").Value = 3.1
$ws.Range("N218").Value = 2
$ws.Range("O218").Value = 3.25
$ws.Range("P218").Value = 3.1
$ws.Range("Q218").Value = 2.05
$ws.Range("R218").Value = 0.25
$ws.Range("S218").Value = 1.95
$ws.Range("T218").Value = 1.85
$ws.Range("U218").Value = 2.25
$ws.Range("V218").Value = 1.9
$ws.Range("W218").Value = 1.9
$ws.Range("X218").Value = -1
$ws.Range("Y218").Value = 2.1
$ws.Range("Z218").Value = -1
$ws.Range("AA218").Value = 0.475
$ws.Range("AB218").Value = -0.5
$ws.Range("AC218").Value = -1
$ws.Range("AD218").Value = 0.8999999999999999

# Row 219
$ws.Range("B219").Value = 7738703
$ws.Range("E219").Value = 'FC Anyang'
$ws.Range("F219").Value = 'Gyeongnam FC'
$ws.Range("G219").Value = 1
$ws.Range("H219").Value = 0
$ws.Range("I219").Value = 1
$ws.Range("J219").Value = 0
$ws.Range("K219").Value = 'H'
$ws.Range("L219").Value = 1.909
$ws.Range("M219").Value = 3.2
$ws.Range("N219").Value = 3.6
$ws.Range("O219").Value = 2.2
$ws.Range("P219").Value = 3
$ws.Range("Q219").Value = 3.1
$ws.Range("R219").Value = -0.25
$ws.Range("S219").Value = 1.95
$ws.Range("T219").Value = 1.85
$ws.Range("U219").Value = 2.5
$ws.Range("V219").Value = 1.975
$ws.Range("W219").Value = 1.825
$ws.Range("X219").Value = 1.2
$ws.Range("Y219").Value = -1
$ws.Range("Z219").Value = -1
$ws.Range("AA219").Value = 0.95
$ws.Range("AB219").Value = -1
$ws.Range("AC219").Value = -1
$ws.Range("AD219").Value = 0.825

# Row 226
$ws.Range("B226").Value = 7739331
$ws.Range("E226").Value = 'Chungnam Asan FC'
$ws.Range("F226").Value = 'Suwon Bluewings'
$ws.Range("G226").Value = 1
$ws.Range("H226").Value = 0
$ws.Range("I226").Value = 0
$ws.Range("J226").Value = 0
$ws.Range("K226").Value = 'H'
$ws.Range("L226").Value = 3.7
$ws.Range("M226").Value = 3.25
$ws.Range("N226").Value = 1.95
$ws.Range("O226").Value = 4.1
$ws.Range("P226").Value = 3.4
$ws.Range("Q226").Value = 1.833
$ws.Range("R226").Value = 0.5
$ws.Range("S226").Value = 1.95
$ws.Range("T226").Value = 1.85
$ws.Range("U226").Value = 2.5
$ws.Range("V226").Value = 1.925
$ws.Range("W226").Value = 1.875
$ws.Range("X226").Value = 3.1
$ws.Range("Y226").Value = -1
$ws.Range("Z226").Value = -1
$ws.Range("AA226").Value = 0.95
$ws.Range("AB226").Value = -1
$ws.Range("AC226").Value = -1
$ws.Range("AD226").Value = 0.875

# Row 227
$ws.Range("B227").Value = 7738672
$ws.Range("E227").Value = 'Seoul ELand FC'
$ws.Range("F227").Value = 'Bucheon'
$ws.Range("G227").Value = 1
$ws.Range("H227").Value = 1
$ws.Range("I227").Value = 1
$ws.Range("J227").Value = 0
$ws.Range("K227").Value = 'D'
$ws.Range("L227").Value = 2.15
$ws.Range("M227").Value = 3.1
$ws.Range("N227").Value = 3.4
$ws.Range("O227").Value = 2
$ws.Range("P227").Value = 3.1
$ws.Range("Q227").Value = 3.7
$ws.Range("R227").Value = -0.5
$ws.Range("S227").Value = 2
$ws.Range("T227").Value = 1.8
$ws.Range("U227").Value = 2.25
$ws.Range("V227").Value = 1.875
$ws.Range("W227").Value = 1.925
$ws.Range("X227").Value = -1
$ws.Range("Y227").Value = 2.1
$ws.Range("Z227").Value = -1
$ws.Range("AA227").Value = -1
$ws.Range("AB227").Value = 0.8
$ws.Range("AC227").Value = -0.5
$ws.Range("AD227").Value = 0.4625

# Row 228
$ws.Range("B228").Value = 7738706
$ws.Range("E228").Value = 'Jeonnam Dragons'
$ws.Range("F228").Value = 'Chungbuk Cheongju'
$ws.Range("G228").Value = 1
$ws.Range("H228").Value = 1
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = 0
$ws.Range("K228").Value = 'D'
$ws.Range("L228").Value = 1.95
$ws.Range("M228").Value = 3.25
$ws.Range("N228").Value = 3.8
$ws.Range("O228").Value = 1.909
$ws.Range("P228").Value = 3.1
$ws.Range("Q228").Value = 4.333
$ws.Range("R228").Value = -0.5
$ws.Range("S228").Value = 1.95
$ws.Range("T228").Value = 1.85
$ws.Range("U228").Value = 2
$ws.Range("V228").Value = 1.8
$ws.Range("W228").Value = 2
$ws.Range("X228").Value = -1
$ws.Range("Y228").Value = 2.1
$ws.Range("Z228").Value = -1
$ws.Range("AA228").Value = -1
$ws.Range("AB228").Value = 0.8500000000000001
$ws.Range("AC228").Value = 0
$ws.Range("AD228").Value = 0

# Row 239
$ws.Range("O239").Value = 1.95
$ws.Range("Q239").Value = 3.6
$ws.Range("R239").Value = -0.5
$ws.Range("S239").Value = 2
$ws.Range("T239").Value = 1.8

# Row 240
$ws.Range("V240").Value = 2
$ws.Range("W240").Value = 1.8
